# test P7 with -10 percent
# Update computed result values across the scenario-output sheets
# (general/x/U/TBar/Q/L) to reflect the re-run with the new parameter.
$wb = $excel.ActiveWorkbook

# --- Sheet "general" ---
$ws = $wb.Worksheets.Item("general")
$ws.Range("B3").Value = 348.247273367219
$ws.Range("B4").Value = 0.01400017738342285
$ws.Range("B6").Value = 44.59727336721894
$ws.Range("B10").Value = 303.65

# --- Sheet "x" ---
$ws = $wb.Worksheets.Item("x")
$ws.Range("B2").Value = 4
$ws.Range("B3").Value = 6
$ws.Range("B4").Value = 10
$ws.Range("B5").Value = 12
$ws.Range("B6").Value = 3
$ws.Range("B7").Value = 13
$ws.Range("B8").Value = 8
$ws.Range("B9").Value = 11
$ws.Range("B10").Value = 9
$ws.Range("B11").Value = 2
$ws.Range("B12").Value = 1
$ws.Range("B13").Value = 5
$ws.Range("B14").Value = 7

# --- Sheet "U" ---
$ws = $wb.Worksheets.Item("U")
$ws.Range("B5").Value = 2
$ws.Range("B6").Value = 3
$ws.Range("B8").Value = 2
$ws.Range("B9").Value = 2
$ws.Range("B11").Value = 2

# --- Sheet "TBar" ---
$ws = $wb.Worksheets.Item("TBar")
$ws.Range("B4").Value = 20
$ws.Range("B5").Value = 24.04101472405137
$ws.Range("B6").Value = 10
$ws.Range("B7").Value = 27.5860495735166
$ws.Range("B8").Value = 20
$ws.Range("B9").Value = 12.01159140980468
$ws.Range("B10").Value = 12.31224998648502
$ws.Range("B11").Value = 25.63617778285959
$ws.Range("B12").Value = 10
$ws.Range("B13").Value = 28.77116560159149
$ws.Range("B14").Value = 25.33666562565053
$ws.Range("B15").Value = 26.35974282367841

# --- Sheet "Q" ---
$ws = $wb.Worksheets.Item("Q")
$ws.Range("C7").Value = 240.8
$ws.Range("C8").Value = 260.6199999999997
$ws.Range("C9").Value = 235.9049999999997
$ws.Range("C10").Value = 248.4
$ws.Range("C11").Value = 238.7299999999997
$ws.Range("C12").Value = 207.4550000000009
$ws.Range("C13").Value = 208.4800000000009
$ws.Range("C14").Value = 211.9650000000009
$ws.Range("C15").Value = 210.7450000000009
$ws.Range("C16").Value = 214.1950000000008
$ws.Range("C17").Value = 254.3750000000003
$ws.Range("C18").Value = 244.5850000000003
$ws.Range("C19").Value = 246.8800000000003
$ws.Range("C20").Value = 248.2350000000003
$ws.Range("C21").Value = 253.6300000000003
$ws.Range("C22").Value = 66.72999999999949
$ws.Range("C23").Value = 67.4749999999995
$ws.Range("C24").Value = 68.55
$ws.Range("C25").Value = 69.7249999999995
$ws.Range("C26").Value = 67.3949999999995
$ws.Range("C27").Value = 262.7450000000006
$ws.Range("C28").Value = 276.5250000000005
$ws.Range("C29").Value = 245.9150000000006
$ws.Range("C30").Value = 272.2100000000006
$ws.Range("C31").Value = 252.4
$ws.Range("C32").Value = 107.3799999999999
$ws.Range("C33").Value = 112.2399999999999
$ws.Range("C34").Value = 93.78999999999985
$ws.Range("C35").Value = 108.8349999999998
$ws.Range("C36").Value = 94.77999999999986
$ws.Range("C37").Value = 125.9800000000001
$ws.Range("C38").Value = 127.8
$ws.Range("C39").Value = 127.25
$ws.Range("C40").Value = 130.3100000000002
$ws.Range("C41").Value = 126.85
$ws.Range("C42").Value = 130.8549999999989
$ws.Range("C43").Value = 143.0299999999989
$ws.Range("C44").Value = 122.6349999999989
$ws.Range("C45").Value = 129.2449999999989
$ws.Range("C46").Value = 122.9249999999989
$ws.Range("C47").Value = 192.8650000000004
$ws.Range("C48").Value = 200.2050000000003
$ws.Range("C49").Value = 173.5050000000004
$ws.Range("C50").Value = 199.5650000000004
$ws.Range("C51").Value = 179.5100000000004
$ws.Range("C52").Value = 53.66499999999927
$ws.Range("C53").Value = 51.98999999999927
$ws.Range("C54").Value = 57.97499999999927
$ws.Range("C55").Value = 55.35
$ws.Range("C56").Value = 50.35499999999927
$ws.Range("C57").Value = 262.7450000000006
$ws.Range("C58").Value = 276.5250000000005
$ws.Range("C59").Value = 245.9150000000006
$ws.Range("C60").Value = 272.2100000000006
$ws.Range("C61").Value = 252.4
$ws.Range("C62").Value = 240.8
$ws.Range("C63").Value = 260.6199999999997
$ws.Range("C64").Value = 235.9049999999997
$ws.Range("C65").Value = 248.4
$ws.Range("C66").Value = 238.7299999999997
$ws.Range("C67").Value = 254.3750000000003
$ws.Range("C68").Value = 244.5850000000003
$ws.Range("C69").Value = 246.8800000000003
$ws.Range("C70").Value = 248.2350000000003
$ws.Range("C71").Value = 253.6300000000003

# --- Sheet "L" ---
$ws = $wb.Worksheets.Item("L")
$ws.Range("C17").Value = 5.9
$ws.Range("C18").Value = 12.58
$ws.Range("C19").Value = 13.77
$ws.Range("C20").Value = 14.23
$ws.Range("C21").Value = 13.42
$ws.Range("C22").Value = 0
$ws.Range("C23").Value = 0
$ws.Range("C24").Value = 0
$ws.Range("C25").Value = 0
$ws.Range("C26").Value = 0
$ws.Range("C32").Value = 10.76
$ws.Range("C33").Value = 8.91
$ws.Range("C34").Value = 8.699999999999999
$ws.Range("C35").Value = 14.81
$ws.Range("C36").Value = 5.36
$ws.Range("C37").Value = 3.8
$ws.Range("C38").Value = 3.605
$ws.Range("C39").Value = 5.735
$ws.Range("C40").Value = 4.25
$ws.Range("C41").Value = 3.4
$ws.Range("C47").Value = 4.285
$ws.Range("C48").Value = 6.69
$ws.Range("C49").Value = 3.755
$ws.Range("C50").Value = 5.305
$ws.Range("C51").Value = 2.56
